$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.828.30'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = '2.193.97'
$ws.Range('E3').Value = '  -2.01%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = '294.03'
$ws.Range('E5').Value = '  -3.76%  '
$ws.Range('E6').Value = '  -3.68%  '
$ws.Range('D7').Value = '0.570'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = '0.479'
$ws.Range('E9').Value = '  -7.67%  '
$ws.Range('D10').Value = '32.28'
$ws.Range('E10').Value = '  -6.07%  '
$ws.Range('D11').Value = '0.0771'
$ws.Range('E11').Value = '  -4.50%  '
$ws.Range('E12').Value = '  -1.93%  '
$ws.Range('D13').Value = '6.71'
$ws.Range('E13').Value = '  -5.62%  '
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('D15').Value = '2.272.38'
$ws.Range('E15').Value = '  -3.85%  '
$ws.Range('D16').Value = '13.03'
$ws.Range('E16').Value = '  -3.46%  '
$ws.Range('D17').Value = '0.761'
$ws.Range('E17').Value = '  -8.36%  '
$ws.Range('D18').Value = '43.604.84'
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('D19').Value = '0.0₃0878'
$ws.Range('E19').Value = '  -8.30%  '
$ws.Range('E20').Value = '  -7.98%  '
$ws.Range('E21').Value = '  -12.72%  '
$ws.Range('D22').Value = '62.94'
$ws.Range('E22').Value = '  -4.03%  '
$ws.Range('D23').Value = '228.36'
$ws.Range('E23').Value = '  -3.49%  '
$ws.Range('E24').Value = '  -11.88%  '
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('E26').Value = '  -7.88%  '
$ws.Range('D27').Value = '2.19'
$ws.Range('E27').Value = '  -0.37%  '
$ws.Range('D28').Value = '35.42'
$ws.Range('E28').Value = '  -8.18%  '
$ws.Range('D29').Value = '9.13'
$ws.Range('E29').Value = '  -6.20%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '149.11'
$ws.Range('E30').Value = '  -2.78%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '18.89'
$ws.Range('E31').Value = '  -5.26%  '
$ws.Range('E32').Value = '  -10.44%  '
$ws.Range('D34').Value = '0.0726'
$ws.Range('E34').Value = '  -8.38%  '
$ws.Range('D35').Value = '0.115'
$ws.Range('E35').Value = '  -2.84%  '
$ws.Range('D36').Value = '2.87'
$ws.Range('E36').Value = '  -7.20%  '
$ws.Range('D37').Value = '0.100'
$ws.Range('E37').Value = '  -7.22%  '
$ws.Range('D38').Value = '1.63'
$ws.Range('E38').Value = '  -7.22%  '
$ws.Range('D39').Value = '13.27'
$ws.Range('E39').Value = '  -9.07%  '
$ws.Range('E40').Value = '  -6.79%  '
$ws.Range('E41').Value = '  -11.06%  '
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('D43').Value = '3.46'
$ws.Range('E43').Value = '  -8.77%  '
$ws.Range('D44').Value = '1.736.56'
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('D45').Value = '1.62'
$ws.Range('E45').Value = '  +0.88%  '
$ws.Range('D46').Value = '67.70'
$ws.Range('E46').Value = '  -0.29%  '
$ws.Range('D47').Value = '73.31'
$ws.Range('E47').Value = '  -8.41%  '
$ws.Range('D48').Value = '0.170'
$ws.Range('E48').Value = '  -10.33%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.414.15'
$ws.Range('E49').Value = '  -2.07%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '90.83'
$ws.Range('E50').Value = '  -8.31%  '
$ws.Range('D51').Value = '7.39'
$ws.Range('E51').Value = '  -8.98%  '
